# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) to the values from the latest scrape.
#
# Price cells that look like plain numbers (single '.') are written with a
# leading apostrophe so Excel keeps them as text (matching the original
# t="inlineStr" cells) instead of auto-converting to a number, and the
# style is reset to "Normal" right after so no quote-prefix formatting
# lingers on the cell. Price cells that already contain two dots (e.g.
# "26.460.63") are never auto-parsed as numbers, so no extra handling is
# needed for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.460.63"
$ws.Range("E2").Value = "  -3.31%  "

$ws.Range("D3").Value = "1.804.43"
$ws.Range("E3").Value = "  -3.02%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").Value = "'1.006"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").Value = "'308.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.30%  "

$ws.Range("D7").Value = "'0.4534"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.81%  "

$ws.Range("D8").Value = "'0.3640"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.13%  "

$ws.Range("D9").Value = "'0.07105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.77%  "

$ws.Range("D10").Value = "'0.8682"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.66%  "

$ws.Range("D11").Value = "'0.07784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("D12").Value = "'19.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.83%  "

$ws.Range("D13").Value = "1.804.37"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").Value = "'5.265"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").Value = "'6.308"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.81%  "

$ws.Range("D16").Value = "'86.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.72%  "

$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "'0.000008541"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.40%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").Value = "26.512.49"
$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").Value = "'14.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.15%  "

$ws.Range("D22").Value = "'4.956"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.54%  "

$ws.Range("D23").Value = "2.021.31"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Value = "'10.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").Value = "'1.974"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.15%  "

$ws.Range("D26").Value = "'150.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "

$ws.Range("D27").Value = "'17.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").Value = "'1.984"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.37%  "

$ws.Range("D29").Value = "'112.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.84%  "

$ws.Range("D30").Value = "'4.859"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.61%  "

$ws.Range("D31").Value = "'0.08634"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.35%  "

$ws.Range("D32").Value = "'3.024"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").Value = "'0.7242"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.27%  "

$ws.Range("D34").Value = "'4.423"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").Value = "'1.113"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.26%  "

$ws.Range("D36").Value = "'2.489"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.06%  "

$ws.Range("D37").Value = "'1.073"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.63%  "

$ws.Range("D38").Value = "'0.01904"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "

$ws.Range("D39").Value = "'0.05068"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.71%  "

$ws.Range("D40").Value = "'2.864"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.57%  "

$ws.Range("D41").Value = "'6.895"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("D42").Value = "'0.4903"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("D43").Value = "'0.1565"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.91%  "

$ws.Range("E44").Value = "  -3.74%  "

$ws.Range("D45").Value = "'1.006"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").Value = "'0.4590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.67%  "

$ws.Range("D47").Value = "'9.960"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'101.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("D49").Value = "'1.576"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.34%  "

$ws.Range("D50").Value = "'0.05990"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.68%  "

$ws.Range("D51").Value = "'63.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.15%  "
